# Generate Report for Handoff
# Adds a new row (row 3) to each of the three worksheets (Overview, zh-cn,
# de-de) for the file 7884b7ee-a1e1-4556-ae56-6cb857d963d9.md, which has
# just been marked "Ready for handoff" / "True" for localization, mirroring
# the existing row 2 (b288ff8f-0b98-48a4-bae0-bda0a8feda01.md) that is
# already "In Translation".

$wb = $excel.ActiveWorkbook

$newFile     = "7884b7ee-a1e1-4556-ae56-6cb857d963d9.md"
$newPath     = "e2e\7884b7ee-a1e1-4556-ae56-6cb857d963d9.md"
$commit      = "f3854b38cfac50751c213b07cc6d35b891a2f475"
$baseUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/" + $commit + "/e2e/" + $newFile
$dateFormat  = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $baseUrl, "", "", $newPath)

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2017-01-03 07:19:48"
$wsOverview.Range("G3").NumberFormat = $dateFormat

$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $baseUrl, "", "", $newFile)

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "7884b7ee-a1e1-4556-ae56-6cb857d963d9.cebea20dc7eed3c0dd42d3691864333d4588c5af.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2017-01-03 07:19:38"
$wsZhCn.Range("H3").NumberFormat = $dateFormat
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = ""
$wsZhCn.Range("L3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L3").NumberFormat = $dateFormat
$wsZhCn.Range("M3").Value = ""
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "True"
$wsZhCn.Range("P3").Value = ""
$wsZhCn.Range("Q3").Value = "False"
$wsZhCn.Range("R3").Value = ""

$wsZhCn.Columns.Item(3).ColumnWidth = 16.33

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:R3"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $baseUrl, "", "", $newFile)

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "7884b7ee-a1e1-4556-ae56-6cb857d963d9.cebea20dc7eed3c0dd42d3691864333d4588c5af.de-de.xlf"
$wsDeDe.Range("H3").Value = "2017-01-03 07:19:48"
$wsDeDe.Range("H3").NumberFormat = $dateFormat
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = ""
$wsDeDe.Range("L3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L3").NumberFormat = $dateFormat
$wsDeDe.Range("M3").Value = ""
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "True"
$wsDeDe.Range("P3").Value = ""
$wsDeDe.Range("Q3").Value = "False"
$wsDeDe.Range("R3").Value = ""

$wsDeDe.Columns.Item(3).ColumnWidth = 16.33

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:R3"))
